$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A105").Value = "modificacion de oc"
$ws.Range("B105").Value = "no comenzado"

$ws.Range("A106").Value = "ot listado arreglar filtro por cliente"
$ws.Range("B106").Value = "no comenzado"

$ws.Range("C104").Select()
